$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New journal entries (row, date serial, activity (col B), problematique (col C or $null))
$entries = @(
    @(13, 44972, "Montage", $null),
    @(14, 44979, "Montage", $null),
    @(15, 44986, "Montage et tests", $null),
    @(16, 44993, "Programmation firmware", "Tests et lectures"),
    @(17, 45000, "Programmation firmware", "Centrale inertiell BNO055"),
    @(18, 45007, "Programmation firmware", "Centrale inertiell BNO055"),
    @(19, 45014, "Programmation firmware", "Centrale inertiell BNO055"),
    @(20, 45021, "Programmation firmware", "Carte SD"),
    @(21, 45028, "Vacances", $null),
    @(22, 45035, "Vacances", $null),
    @(23, 45042, "Design mécanique et tests", $null),
    @(24, 45049, "Programmation firmware", "Carte SD"),
    @(25, 45056, "Implémentation capteur de pression", $null),
    @(26, 45063, "Programmation firmware", "ADC - Capteur de pression"),
    @(27, 45070, "Programmation firmware", "ADC - Capteur de pression"),
    @(28, 45077, "Programmation firmware ", "Gestion des temps processeur")
)

# Copy column A's date-cell format (border + centered + date numfmt) down the
# new rows so the new date cells reuse the existing style instead of Excel
# minting a duplicate one.
$ws.Range("A12").Copy()
$ws.Range("A13:A28").PasteSpecial(-4122)

# Give every new row the same explicit height as the existing ones.
$ws.Range("A13:A28").RowHeight = 45.75

foreach ($entry in $entries) {
    $row = $entry[0]
    $dateSerial = $entry[1]
    $activite = $entry[2]
    $problematique = $entry[3]

    $ws.Cells.Item($row, 1).Value = $dateSerial
    $ws.Cells.Item($row, 2).Value = $activite

    if ($problematique -ne $null) {
        $ws.Cells.Item($row, 3).Value = $problematique
    }
}

# Move the selection to match the saved view state (C29, one past the last
# data row).
$ws.Range("C29").Select()
